# Apply the "ops" reshuffle to the order tracking sheet:
#  - The "Client" column (previously column E) moves to column G.
#  - "Typist" / "Typist QC" (previously F/G) shift left to E/F.
#  - "Product Name" (previously H) and "Lob" (previously J) swap places.
#  - The active selection on the sheet moves from E12 to F4.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Move column E ("Client") to its new home just before the old column H,
# which lands it on column G after the shift. Excel's Cut + Insert keeps
# the exact custom width of the moved column (12.6640625) attached to it.
$ws.Columns("E").Cut()
$ws.Columns("H").Insert()

# The cut/insert leaves behind placeholder <col> entries (now-empty E:F)
# that Excel marks with customWidth="0"; clear them so no stray column
# width definitions remain for E:F (matching the original template, which
# only ever had one custom-width entry there).
$ws.Columns("E:F").ClearFormats()

# ClearFormats() reset E:F (row 1 header + rows 2-3 body) back to the
# default style, so reapply the correct direct formatting (cellXf) by
# copying it over from the still-intact neighboring column D, which has
# the same header/body styling.
$ws.Range("D1").Copy()
$ws.Range("E1:F1").PasteSpecial(-4122)
$ws.Range("D2").Copy()
$ws.Range("E2:F3").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Swap "Product Name" (H) and "Lob" (J) values for every data row.
$tmpH1 = $ws.Range("H1").Value2
$tmpH2 = $ws.Range("H2").Value2
$tmpH3 = $ws.Range("H3").Value2
$ws.Range("H1").Value = $ws.Range("J1").Value2
$ws.Range("H2").Value = $ws.Range("J2").Value2
$ws.Range("H3").Value = $ws.Range("J3").Value2
$ws.Range("J1").Value = $tmpH1
$ws.Range("J2").Value = $tmpH2
$ws.Range("J3").Value = $tmpH3

# Update the saved cursor/selection position on the sheet.
$ws.Range("F4").Select() | Out-Null
